# edit.ps1 - Applies the "2071 a 2130" ordenanza formatting/content edit.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# STEP 1: Text edits (splits / insertions / deletions) -- do these first
# while paragraph/character indices are still simple, working top-down.
# ---------------------------------------------------------------------

# 1a. Split "VISTO: " away from the rest of that paragraph.
$rng = $d.Content
$rng.Find.Execute("VISTO: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# 1b. New leading run " " before "La Ordenanza".
$rng = $d.Content
$rng.Find.Execute("La Ordenanza N", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$rng.InsertBefore(" ")

# 1c. Split "CONSIDERANDO: " away from the rest of that paragraph.
$rng = $d.Content
$rng.Find.Execute("CONSIDERANDO: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# 1d. New leading run " " before "Que se hace".
$rng = $d.Content
$rng.Find.Execute("Que se hace", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$rng.InsertBefore(" ")

# 1e. Drop the "POR " lead-in before "EL CONCEJO DELIBERANTE...".
$rng = $d.Content
$rng.Find.Execute("POR EL CONCEJO", $false, $false, $false, $false, $false, $true, 1, $false, "EL CONCEJO", 2)

# 1f. Collapse the long run of spaces before each "(Presupuesto..." to a
# single space (both occurrences get replaced by ReplaceAll).
$rng = $d.Content
$rng.Find.Execute("                    (", $false, $false, $false, $false, $false, $true, 1, $false, " (", 2)

# ---------------------------------------------------------------------
# STEP 2: Paragraph + character formatting, now that the document has
# settled into its final 10-paragraph shape.
# ---------------------------------------------------------------------

# Para 1: "Yerba Buena, 23 de Diciembre de 2015" (date line)
$p = $d.Paragraphs(1)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 12

# Para 2: "ORDENANZA Nº 2023"
$p = $d.Paragraphs(2)
$p.Format.KeepWithNext = 1
$p.Format.SpaceBefore = 12
$p.Format.SpaceAfter = 18
$p.Range.Font.Bold = 1

# Para 3: "VISTO: "
$p = $d.Paragraphs(3)
$p.Format.KeepWithNext = 1
$p.Format.SpaceBefore = 12
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0
$p.Range.Font.Bold = 1

# Para 4: " La Ordenanza Nº 2.012 ..."
$p = $d.Paragraphs(4)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

# Para 5: "CONSIDERANDO: "
$p = $d.Paragraphs(5)
$p.Format.KeepWithNext = 1
$p.Format.SpaceBefore = 12
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0
$p.Range.Font.Bold = 1

# Para 6: " Que se hace necesario ..."
$p = $d.Paragraphs(6)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

# Para 7: "EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA"
$p = $d.Paragraphs(7)
$p.Format.KeepWithNext = 1
$p.Format.SpaceBefore = 18
$p.Format.SpaceAfter = 18
$p.Format.LeftIndent = $word.InchesToPoints(1984 / 1440)
$p.Format.RightIndent = $word.InchesToPoints(1984 / 1440)
$p.Format.Alignment = 1
$p.Range.Font.Bold = 1

# Para 8: "ARTICULO PRIMERO: PRORRÓGASEla vigencia ..."
$p = $d.Paragraphs(8)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

# Para 9: "ARTICULO SEGUNDO: PRORROGASE la vigencia ..."
$p = $d.Paragraphs(9)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

# Para 10: "ARTICULO TERCERO: COMUNÍQUESE, ..."
$p = $d.Paragraphs(10)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

# ---------------------------------------------------------------------
# STEP 3: Underline the "ARTICULO ###:" labels in the three operative
# articles.
# ---------------------------------------------------------------------
$labels = "ARTICULO PRIMERO:", "ARTICULO SEGUNDO:", "ARTICULO TERCERO:"
foreach ($label in $labels) {
    $rng = $d.Content
    $rng.Find.Execute($label, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng.Font.Underline = 1
}

# ---------------------------------------------------------------------
# STEP 5: Section properties - restart page numbering at 2977.
# ---------------------------------------------------------------------
$sec = $d.Sections(1)
$hdr = $sec.Headers.Item(1)
$pns = $hdr.PageNumbers
$pns.RestartNumberingAtSection = 1
$pns.StartingNumber = 2977

foreach ($p in $d.Paragraphs) {
    Write-Host "PARA:" $p.Range.Text
}
